$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, border, centered) from H1 into the new
# header cells I1 and J1, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column values for rows 2-24 (I0, IF)
$data = @(
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(3, 4),
    @(7, 7),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(7, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(3, 3),
    @(8, 8),
    @(7, 7),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

Write-Host "edit applied"
